# Append the next day's (2025-10-29, Excel serial 45959) readings for both
# charging stations as two new rows at the bottom of the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = 45959
$firstDataCol = 3   # column C (hourly buckets start here; A = date, B = station)
$lastDataCol  = 26  # column Z

$stationRows = @(
    @{
        Row     = 118
        Station = "四方坪站充电量(kw)"
        Values  = @(998.51400000000001, 1190.1950000000002, 191.11199999999999, 384.40000000000009, 260.87400000000002, 720.13000000000011, 174.32399999999998, 188.15299999999999, 191.16199999999998, 158.39099999999999, 205.76999999999998, 288.36, 718.14800000000002, 1420.376, 426.06999999999994, 430.28800000000007, 179.30999999999997, 244.65500000000003, 61.646000000000001, 85.034000000000006, 78.62299999999999, 32.893999999999998, 63.78, 9.08)
    },
    @{
        Row     = 119
        Station = "高岭站充电量(kw)"
        Values  = @(311.92400000000004, 309.21100000000001, 161.53100000000001, 62.046999999999997, 112.52000000000001, 63.858999999999995, 36.969000000000001, 91.103000000000009, 222.23599999999999, 170.76900000000001, 169.24199999999999, 145.19899999999998, 439.50899999999996, 402.39000000000004, 224.57100000000003, 212.78, 66.461000000000013, 31.405000000000001, 57.529000000000003, 23.407, 27.57, 14.055, 7.0629999999999997, 0)
    }
)

foreach ($entry in $stationRows) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $newDate        # column A: 日期
    $ws.Cells.Item($r, 2).Value = $entry.Station   # column B: 站点

    for ($col = $firstDataCol; $col -le $lastDataCol; $col++) {
        $ws.Cells.Item($r, $col).Value = $entry.Values[$col - $firstDataCol]
    }
}

# Match the author's final cursor position after the append.
$ws.Range("N122").Select()
